# Apply updates described by the commit: "Updated capital structure database"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (company #2): updated metric values ---
$ws.Range("B2").Value2 = "'2"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Value2 = [double]"-0.015675"
$ws.Range("E2").Value2 = [double]"-0.0513"
$ws.Range("G2").Value2 = [double]"0"
$ws.Range("K2").Value2 = [double]"13.16"
$ws.Range("L2").Value2 = [double]"0.03440522875816993"
$ws.Range("M2").Value2 = [double]"17.973"
$ws.Range("N2").Value2 = [double]"0.0610703363914373"
$ws.Range("O2").Value2 = [double]"1.365729483282675"
$ws.Range("P2").Value2 = [double]"17.973"
$ws.Range("Q2").Value2 = [double]"0.0610703363914373"
$ws.Range("R2").Value2 = [double]"1.365729483282675"
$ws.Range("U2").Value2 = [double]"0.006"
$ws.Range("V2").Value2 = [double]"2.038735983690112e-05"
$ws.Range("W2").Value2 = [double]"0.06621547029856485"
$ws.Range("X2").Value2 = [double]"0.05501979681991646"
$ws.Range("Y2").Value2 = [double]"0.01119567347864839"
$ws.Range("Z2").Value2 = [double]"0.2838813294403864"
$ws.Range("AB2").Value2 = [double]"0.03523317733323548"
$ws.Range("AC2").Value2 = [double]"-0.03523317733323548"
$ws.Range("AD2").Value2 = [double]"1163.5"
$ws.Range("AF2").Value2 = [double]"1163.5"
$ws.Range("AG2").Value2 = [double]"1163.494"
$ws.Range("AH2").Value2 = [double]"0.7981204554808616"
$ws.Range("AI2").Value2 = [double]"0.8569639832068939"
$ws.Range("AJ2").Value2 = [double]"0.798119624583446"
$ws.Range("AK2").Value2 = [double]"0.8569633510938399"

# --- Row 3 (Salafin S.A.): updated metric values ---
$ws.Range("D3").Value2 = [double]"-0.0299"
$ws.Range("E3").Value2 = [double]"-0.0741"
$ws.Range("K3").Value2 = [double]"7.86"
$ws.Range("L3").Value2 = [double]"0.1070844686648501"
$ws.Range("M3").Value2 = [double]"13.393"
$ws.Range("N3").Value2 = [double]"0.08423270440251572"
$ws.Range("O3").Value2 = [double]"1.703944020356234"
$ws.Range("P3").Value2 = [double]"13.393"
$ws.Range("Q3").Value2 = [double]"0.08423270440251572"
$ws.Range("R3").Value2 = [double]"1.703944020356234"
$ws.Range("V3").Value2 = [double]"3.773584905660378e-05"
$ws.Range("W3").Value2 = [double]"0.07507163323782234"
$ws.Range("X3").Value2 = [double]"0.03005133547844983"
$ws.Range("Y3").Value2 = [double]"0.04502029775937252"
$ws.Range("Z3").Value2 = [double]"0.2208887310634559"
$ws.Range("AB3").Value2 = [double]"0.02922580070869248"
$ws.Range("AC3").Value2 = [double]"-0.02922580070869248"
$ws.Range("AD3").Value2 = [double]"177.4"
$ws.Range("AF3").Value2 = [double]"177.4"
$ws.Range("AG3").Value2 = [double]"177.394"
$ws.Range("AH3").Value2 = [double]"0.5273483947681332"
$ws.Range("AI3").Value2 = [double]"0.6446220930232559"
$ws.Range("AJ3").Value2 = [double]"0.5273399644464527"
$ws.Range("AK3").Value2 = [double]"0.644614344789494"

# --- Row 4 (renamed to Maroc Leasing S.A.): updated metric values ---
$ws.Range("B4").Value2 = "Maroc Leasing S.A. (CBSE:MLE)"
$ws.Range("D4").Value2 = [double]"-0.00145"
$ws.Range("E4").Value2 = [double]"-0.0285"
$ws.Range("G4").Value2 = [double]"0"
$ws.Range("K4").Value2 = [double]"5.3"
$ws.Range("L4").Value2 = [double]"0.01714655451310255"
$ws.Range("M4").Value2 = [double]"4.58"
$ws.Range("N4").Value2 = [double]"0.03385070214338507"
$ws.Range("O4").Value2 = [double]"0.8641509433962264"
$ws.Range("P4").Value2 = [double]"4.58"
$ws.Range("Q4").Value2 = [double]"0.03385070214338507"
$ws.Range("R4").Value2 = [double]"0.8641509433962264"
$ws.Range("U4").Value2 = [double]"0"
$ws.Range("V4").Value2 = [double]"0"
$ws.Range("W4").Value2 = [double]"0.05735930735930735"
$ws.Range("X4").Value2 = [double]"0.07998825816138311"
$ws.Range("Y4").Value2 = [double]"-0.02262895080207575"
$ws.Range("Z4").Value2 = [double]"0.3045020195054675"
$ws.Range("AB4").Value2 = [double]"0.04124055395777848"
$ws.Range("AC4").Value2 = [double]"-0.04124055395777848"
$ws.Range("AD4").Value2 = [double]"986.1"
$ws.Range("AF4").Value2 = [double]"986.1"
$ws.Range("AG4").Value2 = [double]"986.1"
$ws.Range("AH4").Value2 = [double]"0.8793472445157838"
$ws.Range("AI4").Value2 = [double]"0.9109468822170901"
$ws.Range("AJ4").Value2 = [double]"0.8793472445157838"
$ws.Range("AK4").Value2 = [double]"0.9109468822170901"

# --- Row 5 (Maroc Leasing S.A. old entry) is removed entirely; data now lives in row 4 ---
$ws.Rows(5).Delete()

